$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("J2").Value = 1.06
$ws.Range("K2").Value = 10

# Row 9
$ws.Range("G9").Value = 1.95
$ws.Range("H9").Value = 3.7
$ws.Range("L9").Value = 1.19
$ws.Range("Z9").Value = 14
$ws.Range("AA9").Value = 7.4

# Row 11
$ws.Range("N11").Value = 1.44
$ws.Range("O11").Value = 2.42
$ws.Range("AD11").Value = 200
$ws.Range("AE11").Value = 10

# Row 20
$ws.Range("N20").Value = 2.18
$ws.Range("AE20").Value = 11
$ws.Range("AJ20").Value = 75

# Row 32
$ws.Range("G32").Value = 2.38
$ws.Range("I32").Value = 2.88
$ws.Range("AE32").Value = 10
$ws.Range("AJ32").Value = 34

# Row 37
$ws.Range("K37").Value = 13

# Row 39
$ws.Range("G39").Value = 1.82
$ws.Range("H39").Value = 2.92
$ws.Range("I39").Value = 5
$ws.Range("J39").Value = 1.15
$ws.Range("K39").Value = 4.2
$ws.Range("L39").Value = 1.65
$ws.Range("M39").Value = 2
$ws.Range("N39").Value = 2.82
$ws.Range("O39").Value = 1.32
$ws.Range("P39").Value = 1.65
$ws.Range("Q39").Value = 2
$ws.Range("R39").Value = 2.5
$ws.Range("S39").Value = 1.4
$ws.Range("T39").Value = 4.35
$ws.Range("V39").Value = 10
$ws.Range("W39").Value = 14.5
$ws.Range("X39").Value = 22
$ws.Range("Y39").Value = 60
$ws.Range("Z39").Value = 4.5
$ws.Range("AA39").Value = 6.5
$ws.Range("AB39").Value = 28
$ws.Range("AC39").Value = 250
$ws.Range("AE39").Value = 8.5
$ws.Range("AF39").Value = 26
$ws.Range("AG39").Value = 20
$ws.Range("AI39").Value = 90
$ws.Range("AJ39").Value = 120
